$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates ---
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"
$ws.Range("G10").Value = ""

# --- Numeric cell updates (line item pricing / totals) ---
$ws.Range("C8").Value = 16974.85
$ws.Range("H16").Value = 648.53
$ws.Range("H17").Value = 648.53
$ws.Range("H18").Value = 648.53
$ws.Range("H19").Value = 648.53
$ws.Range("H20").Value = 350.53
$ws.Range("H21").Value = 63.44
$ws.Range("H22").Value = 82.8
$ws.Range("H23").Value = 858.75
$ws.Range("H24").Value = 110.74
$ws.Range("H25").Value = 79.34999999999999
$ws.Range("H26").Value = 478.55
$ws.Range("H27").Value = 187.26
$ws.Range("H28").Value = 73.73
$ws.Range("H29").Value = 18.51
$ws.Range("H30").Value = 285.45
$ws.Range("F31").Value = 3
$ws.Range("H31").Value = 285.45
$ws.Range("H32").Value = 35.58
$ws.Range("H33").Value = 31.08
$ws.Range("H34").Value = 17.2
$ws.Range("H35").Value = 79.56
$ws.Range("H36").Value = 350.53
$ws.Range("H37").Value = 282.51
$ws.Range("H38").Value = 94.17
$ws.Range("H39").Value = 188.34
$ws.Range("H40").Value = 478.55
$ws.Range("H41").Value = 55.18
$ws.Range("H42").Value = 61.83
$ws.Range("H43").Value = 93.23999999999999
$ws.Range("H44").Value = 31.08
$ws.Range("H45").Value = 62.16
$ws.Range("H46").Value = 198.88
$ws.Range("H47").Value = 17.2
$ws.Range("H48").Value = 350.53
$ws.Range("H49").Value = 95.16
$ws.Range("H50").Value = 82.8
$ws.Range("H51").Value = 110.74
$ws.Range("H52").Value = 376.68
$ws.Range("H53").Value = 94.17
$ws.Range("H54").Value = 188.34
$ws.Range("H55").Value = 478.55
$ws.Range("H56").Value = 55.18
$ws.Range("H57").Value = 187.26
$ws.Range("H58").Value = 61.83
$ws.Range("H59").Value = 17.37
$ws.Range("H60").Value = 18.51
$ws.Range("H61").Value = 124.32
$ws.Range("H62").Value = 31.08
$ws.Range("H63").Value = 62.16
$ws.Range("H64").Value = 17.2
$ws.Range("H65").Value = 79.56
$ws.Range("H66").Value = 350.53
$ws.Range("H67").Value = 31.72
$ws.Range("H68").Value = 27.6
$ws.Range("H69").Value = 286.25
$ws.Range("H70").Value = 223
$ws.Range("H71").Value = 110.74
$ws.Range("H72").Value = 282.51
$ws.Range("H73").Value = 94.17
$ws.Range("H74").Value = 188.34
$ws.Range("H75").Value = 478.55
$ws.Range("H76").Value = 110.36
$ws.Range("H77").Value = 62.42
$ws.Range("H78").Value = 159
$ws.Range("H79").Value = 61.83
$ws.Range("H80").Value = 17.37
$ws.Range("H81").Value = 6.17
$ws.Range("H82").Value = 95.15000000000001
$ws.Range("H83").Value = 116
$ws.Range("H84").Value = 93.23999999999999
$ws.Range("H85").Value = 31.08
$ws.Range("H86").Value = 62.16
$ws.Range("H87").Value = 198.88
$ws.Range("H88").Value = 34.4
$ws.Range("H89").Value = 26.52
$ws.Range("H90").Value = 478.55
$ws.Range("H91").Value = 478.55
$ws.Range("H92").Value = 478.55
$ws.Range("H93").Value = 478.55
$ws.Range("H94").Value = 478.55
$ws.Range("H95").Value = 15517.95
$ws.Range("H100").Value = 216.17
$ws.Range("H101").Value = 234
$ws.Range("H102").Value = 450.17
$ws.Range("H107").Value = 216.17
$ws.Range("H108").Value = 216.17
$ws.Range("H113").Value = 790.5599999999999
$ws.Range("H114").Value = 790.5599999999999

Write-Host "Applied all changes"
